# Updating the country bounding box mapping/expansion routines
#
# - Shorten the two IFRC GO taxonomy labels used in column C:
#     "GO-Appeal"        -> "GO-App"   (rows 31-33)
#     "GO-Field Reports" -> "GO-FR"    (rows 34-55)
# - C32/C33 had picked up the wrong (header-ish) cell style; line them back
#   up with the rest of column C (same look as C31/C34).
# - Leave the view zoomed/scrolled/selected where the editor left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the taxonomy labels (every cell sharing the string updates together)
$ws.Range("C31:C33").Replace("GO-Appeal", "GO-App")
$ws.Range("C34:C55").Replace("GO-Field Reports", "GO-FR")

# Re-align the formatting on C32/C33 with the rest of the column
$ws.Range("C31").Copy()
$ws.Range("C32:C33").PasteSpecial(-4122)

# Update the view: zoom out a bit, scroll up, and move the active cell
$excel.ActiveWindow.Zoom = 110
$excel.ActiveWindow.ScrollRow = 20
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F40").Select()
